$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the "Sunny" value from I2 (leftover value causing the exceed leave
# balance validation issue) - clear the cell entirely.
$ws.Range("I2").ClearContents()

# Update the active selection to I2, matching the saved selection state.
$ws.Range("I2").Select()
